$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Indicator text (B4) updated to the new 16.1.4 wording, and given its own
#    distinct font/style (mirrors the extra cellXf added in the real edit).
$ws.Range("B4").Value = "16.1.4 Proportion of population that feel safe walking alone around the area they live after dark"
$ws.Range("B4").Font.Name = "Calibri"

# 2. Data reporter block (B6:B10) - refreshed contact details for the
#    organization / focal point.
$ws.Range("B6").Value = "National Statistical Committee of the Kyrgyz Republic (Department of Household Statistics)"
$ws.Range("B7").Value = "Kalymbetova Yryskan"
$ws.Range("B8").Value = "yryskan.kalymbetova@gmail.com "
$ws.Range("B9").Value = "(0312) 32 46 55"
$ws.Range("B10").Value = "www.stat.gov.kg"
